# Generate Report for Handoff
# The "b.md" entry moves from "Handed back" status to "Ready for handoff":
# a new handoff file/timestamp is recorded for b.md in each language sheet,
# while the already-handed-back handback info stays as-is.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview" - summary row for b.md
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("D3").Value = "2016-29-12 22:29:26"

# ---------------------------------------------------------------------------
# Helper to update a language detail sheet (zh-cn / de-de)
# ---------------------------------------------------------------------------
function Update-LangSheet {
    param(
        [string]$sheetName,
        [string]$newHandoffFile,
        [string]$newHandoffDate,
        [string]$newHandoffUrl
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3 corresponds to "b.md"
    $ws.Range("C3").Value = $newStatus
    $ws.Range("D3").Value = $newHandoffFile
    $ws.Range("E3").Value = $newHandoffDate

    # Refresh the hyperlink on D3 so its cached display text matches the new
    # handoff file name (remove the stale hyperlink entry first, then add the
    # replacement pointing at the new handoff file).
    $hyperlinks = $ws.Hyperlinks
    $target = $null
    for ($i = 1; $i -le $hyperlinks.Count; $i++) {
        $candidate = $hyperlinks.Item($i)
        if ($candidate.Range.Address() -eq '$D$3') {
            $target = $candidate
        }
    }
    if ($target -ne $null) {
        $target.Delete()
    }
    $hyperlinks.Add($ws.Range("D3"), $newHandoffUrl, "", "", $newHandoffFile) | Out-Null
}

Update-LangSheet `
    "zh-cn" `
    "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" `
    "2016-03-12 22:29:23" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc7b3f1d27a08467cbd1b9b8696601f854bea225/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

Update-LangSheet `
    "de-de" `
    "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" `
    "2016-03-12 22:29:26" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a05290431e964e1319928e973b9d26378bc92e87/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
